$wb = $excel.ActiveWorkbook

# --- Update status / timestamp values (Report generated for handoff) ---

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-20 23:02:48"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-20 23:02:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-20 23:02:48"

# --- Widen the "Status" columns to fit the new, longer status text ---

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
